# Apply the authored change:
#  - Update the "C" column (values previously shared-string "dfg", idx 11)
#    on worksheet "A" to a new shared string "new name".
#  - Worksheet "A" becomes the active sheet / selected tab, with the
#    selection moved to C3:C6 (active cell C3).
#  - Worksheet "Q" keeps its own remembered selection, now C2
#    (active cell C2), and loses the "tabSelected" flag since "A" is
#    now the foreground sheet.

$wb = $excel.ActiveWorkbook

$wsA = $wb.Worksheets.Item("A")
$wsQ = $wb.Worksheets.Item("Q")

# Update column C (rows 2-6) on sheet "A" to the new value. This
# registers "new name" as a brand-new shared string, reused by every
# cell in the range.
$wsA.Range("C2:C6").Value = "new name"

# Set sheet "Q"'s remembered selection first (it is not the final active
# sheet, but Excel still records a per-sheet cursor position).
$wsQ.Activate()
$wsQ.Range("C2").Select()

# Finally activate sheet "A" and set its selection - this also makes "A"
# the saved active tab / tabSelected sheet.
$wsA.Activate()
$wsA.Range("C3:C6").Select()
